# Apply the Natmi (Hou advice) re-run update: rows 2-7 change, rows 8-10 are new,
# filling out the full 3x3 ECs/FAPs/sCs sending x target cluster grid.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Hbegf"
$ws.Range("C2").Value = "Erbb2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 13.739149
$ws.Range("H2").Value = 41.217447
$ws.Range("I2").Value = 0.6130043224686931
$ws.Range("J2").Value = 0.6130043224686931
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 1.534538333333333
$ws.Range("N2").Value = 4.603615
$ws.Range("O2").Value = 0.1494637976135089
$ws.Range("P2").Value = 0.1494637976135089
$ws.Range("Q2").Value = 21.08325080787833
$ws.Range("R2").Value = 189.749257270905
$ws.Range("S2").Value = 0.09162195398966692
$ws.Range("T2").Value = 0.09162195398966692

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Hbegf"
$ws.Range("C3").Value = "Erbb2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 13.739149
$ws.Range("H3").Value = 41.217447
$ws.Range("I3").Value = 0.6130043224686931
$ws.Range("J3").Value = 0.6130043224686931
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 4.973328333333334
$ws.Range("N3").Value = 14.919985
$ws.Range("O3").Value = 0.4844014146353658
$ws.Range("P3").Value = 0.4844014146353658
$ws.Range("Q3").Value = 68.32929899758834
$ws.Range("R3").Value = 614.9636909782951
$ws.Range("S3").Value = 0.2969401609814289
$ws.Range("T3").Value = 0.2969401609814289

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Hbegf"
$ws.Range("C4").Value = "Erbb2"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 13.739149
$ws.Range("H4").Value = 41.217447
$ws.Range("I4").Value = 0.6130043224686931
$ws.Range("J4").Value = 0.6130043224686931
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 3.75909
$ws.Range("N4").Value = 11.27727
$ws.Range("O4").Value = 0.3661347877511252
$ws.Range("P4").Value = 0.3661347877511252
$ws.Range("Q4").Value = 51.64669761441
$ws.Range("R4").Value = 464.82027852969
$ws.Range("S4").Value = 0.2244422074975972
$ws.Range("T4").Value = 0.2244422074975972

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Hbegf"
$ws.Range("C5").Value = "Erbb2"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 4.034036666666666
$ws.Range("H5").Value = 12.10211
$ws.Range("I5").Value = 0.1799879973398545
$ws.Range("J5").Value = 0.1799879973398545
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 1.534538333333333
$ws.Range("N5").Value = 4.603615
$ws.Range("O5").Value = 0.1494637976135089
$ws.Range("P5").Value = 0.1494637976135089
$ws.Range("Q5").Value = 6.19038390307222
$ws.Range("R5").Value = 55.71345512764999
$ws.Range("S5").Value = 0.0269016896072648
$ws.Range("T5").Value = 0.0269016896072648

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Hbegf"
$ws.Range("C6").Value = "Erbb2"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 4.034036666666666
$ws.Range("H6").Value = 12.10211
$ws.Range("I6").Value = 0.1799879973398545
$ws.Range("J6").Value = 0.1799879973398545
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 4.973328333333334
$ws.Range("N6").Value = 14.919985
$ws.Range("O6").Value = 0.4844014146353658
$ws.Range("P6").Value = 0.4844014146353658
$ws.Range("Q6").Value = 20.06258885203889
$ws.Range("R6").Value = 180.56329966835
$ws.Range("S6").Value = 0.087186440528812
$ws.Range("T6").Value = 0.08718644052881198

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Hbegf"
$ws.Range("C7").Value = "Erbb2"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 4.034036666666666
$ws.Range("H7").Value = 12.10211
$ws.Range("I7").Value = 0.1799879973398545
$ws.Range("J7").Value = 0.1799879973398545
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 3.75909
$ws.Range("N7").Value = 11.27727
$ws.Range("O7").Value = 0.3661347877511252
$ws.Range("P7").Value = 0.3661347877511252
$ws.Range("Q7").Value = 15.1643068933
$ws.Range("R7").Value = 136.4787620397
$ws.Range("S7").Value = 0.06589986720377772
$ws.Range("T7").Value = 0.06589986720377772

# Row 8
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Hbegf"
$ws.Range("C8").Value = "Erbb2"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 4.639623666666666
$ws.Range("H8").Value = 13.918871
$ws.Range("I8").Value = 0.2070076801914524
$ws.Range("J8").Value = 0.2070076801914524
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 1.534538333333333
$ws.Range("N8").Value = 4.603615
$ws.Range("O8").Value = 0.1494637976135089
$ws.Range("P8").Value = 0.1494637976135089
$ws.Range("Q8").Value = 7.119680368740553
$ws.Range("R8").Value = 64.07712331866499
$ws.Range("S8").Value = 0.03094015401657723
$ws.Range("T8").Value = 0.03094015401657723

# Row 9
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Hbegf"
$ws.Range("C9").Value = "Erbb2"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 4.639623666666666
$ws.Range("H9").Value = 13.918871
$ws.Range("I9").Value = 0.2070076801914524
$ws.Range("J9").Value = 0.2070076801914524
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 4.973328333333334
$ws.Range("N9").Value = 14.919985
$ws.Range("O9").Value = 0.4844014146353658
$ws.Range("P9").Value = 0.4844014146353658
$ws.Range("Q9").Value = 23.07437183743722
$ws.Range("R9").Value = 207.669346536935
$ws.Range("S9").Value = 0.1002748131251249
$ws.Range("T9").Value = 0.1002748131251249

# Row 10
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Hbegf"
$ws.Range("C10").Value = "Erbb2"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 4.639623666666666
$ws.Range("H10").Value = 13.918871
$ws.Range("I10").Value = 0.2070076801914524
$ws.Range("J10").Value = 0.2070076801914524
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 3.75909
$ws.Range("N10").Value = 11.27727
$ws.Range("O10").Value = 0.3661347877511252
$ws.Range("P10").Value = 0.3661347877511252
$ws.Range("Q10").Value = 17.44076292913
$ws.Range("R10").Value = 156.96686636217
$ws.Range("S10").Value = 0.07579271304975023
$ws.Range("T10").Value = 0.07579271304975023
